$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Footer date cell: show date+time instead of just the date
$ws.Range("U3").Formula = "=NOW()"
$ws.Range("U3").NumberFormat = "m/d/yyyy h:mm"
